$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibition listings) — rows keyed by event, columns F (想去人数) and G (最低票价)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 4643
$ws1.Range("G2").Value = 70
$ws1.Range("F3").Value = 2525
$ws1.Range("G3").Value = 60
$ws1.Range("G5").Value = 50
$ws1.Range("F6").Value = 62
$ws1.Range("F8").Value = 234
$ws1.Range("F10").Value = 189
$ws1.Range("F11").Value = 178
$ws1.Range("F12").Value = 1735
$ws1.Range("F14").Value = 3871
$ws1.Range("F15").Value = 36
$ws1.Range("F16").Value = 259

# Sheet "全部类型" (all categories combined) — same events, different row numbers
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 4643
$ws4.Range("G2").Value = 70
$ws4.Range("F3").Value = 2525
$ws4.Range("G3").Value = 60
$ws4.Range("G5").Value = 50
$ws4.Range("F7").Value = 62
$ws4.Range("F10").Value = 234
$ws4.Range("F12").Value = 189
$ws4.Range("F13").Value = 178
$ws4.Range("F16").Value = 1735
$ws4.Range("F18").Value = 3871
$ws4.Range("F19").Value = 36
$ws4.Range("F20").Value = 259
